# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
# The worker/period rows (B16:G29) are re-sorted: ORLANDO ENRIQUE MENDEZ MORENO's
# 7 periods (2107 -> 2101, descending) now come first, followed by
# ALEXANDER HENRIQUE MENDEZ MORENO's 7 periods (2107 -> 2101, descending).
# The underlying set of rows/values is unchanged - only the row order differs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("CC", "1232597678", "ORLANDO ENRIQUE MENDEZ MORENO",   "2107", 29260, 908526),
    @("CC", "1232597678", "ORLANDO ENRIQUE MENDEZ MORENO",   "2106", 35112, 908526),
    @("CC", "1232597678", "ORLANDO ENRIQUE MENDEZ MORENO",   "2105", 35112, 908526),
    @("CC", "1232597678", "ORLANDO ENRIQUE MENDEZ MORENO",   "2104", 35112, 908526),
    @("CC", "1232597678", "ORLANDO ENRIQUE MENDEZ MORENO",   "2103", 35112, 908526),
    @("CC", "1232597678", "ORLANDO ENRIQUE MENDEZ MORENO",   "2102", 35112, 908526),
    @("CC", "1232597678", "ORLANDO ENRIQUE MENDEZ MORENO",   "2101", 35112, 908526),
    @("CC", "1237439959", "ALEXANDER HENRIQUE MENDEZ MORENO", "2107", 29260, 877803),
    @("CC", "1237439959", "ALEXANDER HENRIQUE MENDEZ MORENO", "2106", 35112, 877803),
    @("CC", "1237439959", "ALEXANDER HENRIQUE MENDEZ MORENO", "2105", 35112, 877803),
    @("CC", "1237439959", "ALEXANDER HENRIQUE MENDEZ MORENO", "2104", 35112, 877803),
    @("CC", "1237439959", "ALEXANDER HENRIQUE MENDEZ MORENO", "2103", 35112, 877803),
    @("CC", "1237439959", "ALEXANDER HENRIQUE MENDEZ MORENO", "2102", 35112, 877803),
    @("CC", "1237439959", "ALEXANDER HENRIQUE MENDEZ MORENO", "2101", 35112, 877803)
)

$startRow = 16
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 2).Value = $data[0]   # B - Tipo Doc Trabajador
    $ws.Cells.Item($r, 3).Value = $data[1]   # C - N Doc Trabajador
    $ws.Cells.Item($r, 4).Value = $data[2]   # D - Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $data[3]   # E - Periodo Mora
    $ws.Cells.Item($r, 6).Value = $data[4]   # F - Valor Mora
    $ws.Cells.Item($r, 7).Value = $data[5]   # G - Salario Basico
}
